$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet1"): selection moves to the full column L, and it
# loses tabSelected once the second sheet is activated below.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("L1:L1048576").Select()

# --- Sheet2 ("with_start_date" -> "shuffle_col_order"):
# drop the "start date"/"end date" columns (I:J), which shifts the
# remaining CAGR/ref date/label/comment/source columns left, then add a
# trailing "tags" column with a single "x" value, matching Sheet1's layout.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "shuffle_col_order"
$ws2.Columns("I:J").Delete()

$ws2.Range("N1").Value = "tags"
$ws2.Range("N1").Font.Color = 0
$ws2.Range("N2").Value = "x"

$ws2.Activate()
$ws2.Range("N3").Select()
